$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 42 into the new row 43 (same formatting as the source row)
$ws.Range("A42:F42").Copy()
$ws.Range("A43").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("A42:F42").Copy()
$ws.Range("A43").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Update the new row with the new video's title and YouTube id
$ws.Range("F43").Value = "pXe8MpU7uzk"
$ws.Range("C43").Value = "Live Hummingbird Feeder Cam in Peru"

# Reflect the selection left after the paste/edit operations
[void]$ws.Range("A43:B43").Select()
